# Auto-generated edit script: updates Leve profit/price market-data cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (FFXIV crafting leve profit tracker).
# Values reflect a refreshed market-board data pull (scheduled runner update).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1997.375
$ws.Range("I18").Value = 496.5
$ws.Range("K18").Value = 496.5
$ws.Range("M18").Value = -212.5

$ws.Range("H20").Value = 14487.286
$ws.Range("I20").Value = 1840.75
$ws.Range("J20").Value = 31349.334
$ws.Range("K20").Value = 1840.75
$ws.Range("L20").Value = 31349.334
$ws.Range("M20").Value = -1610.75
$ws.Range("N20").Value = -31809.334

$ws.Range("H21").Value = 1500
$ws.Range("I21").Value = 1500
$ws.Range("K21").Value = 1500
$ws.Range("M21").Value = -1032

$ws.Range("H23").Value = 1500
$ws.Range("I23").Value = 1500
$ws.Range("K23").Value = 1500
$ws.Range("M23").Value = -1266

$ws.Range("H33").Value = 424.43478
$ws.Range("I33").Value = 435.5263
$ws.Range("K33").Value = 435.5263
$ws.Range("M33").Value = -206.5263

$ws.Range("H35").Value = 14487.286
$ws.Range("I35").Value = 1840.75
$ws.Range("J35").Value = 31349.334
$ws.Range("K35").Value = 1840.75
$ws.Range("L35").Value = 31349.334
$ws.Range("M35").Value = -1461.75
$ws.Range("N35").Value = -32107.334

$ws.Range("H62").Value = 3103.3333
$ws.Range("I62").Value = 4005
$ws.Range("K62").Value = 4005
$ws.Range("M62").Value = -3381

$ws.Range("H64").Value = 2027808
$ws.Range("I64").Value = 3080322.5
$ws.Range("J64").Value = 3741.5386
$ws.Range("K64").Value = 3080322.5
$ws.Range("L64").Value = 3741.5386
$ws.Range("M64").Value = -3080074.5
$ws.Range("N64").Value = -4237.5386

$ws.Range("H65").Value = 3103.3333
$ws.Range("I65").Value = 4005
$ws.Range("K65").Value = 20025
$ws.Range("M65").Value = -16905

$ws.Range("H67").Value = 2027808
$ws.Range("I67").Value = 3080322.5
$ws.Range("J67").Value = 3741.5386
$ws.Range("K67").Value = 3080322.5
$ws.Range("L67").Value = 3741.5386
$ws.Range("M67").Value = -3079464.5
$ws.Range("N67").Value = -5457.5386

$ws.Range("H76").Value = 3782.037
$ws.Range("I76").Value = 3811.9443
$ws.Range("J76").Value = 3722.2222
$ws.Range("K76").Value = 3811.9443
$ws.Range("L76").Value = 3722.2222
$ws.Range("M76").Value = -3496.9443
$ws.Range("N76").Value = -4352.2222

$ws.Range("H79").Value = 3782.037
$ws.Range("I79").Value = 3811.9443
$ws.Range("J79").Value = 3722.2222
$ws.Range("K79").Value = 3811.9443
$ws.Range("L79").Value = 3722.2222
$ws.Range("M79").Value = -2719.9443
$ws.Range("N79").Value = -5906.2222

$ws.Range("H98").Value = 1047.5
$ws.Range("I98").Value = 763.75
$ws.Range("J98").Value = 2750
$ws.Range("K98").Value = 763.75
$ws.Range("L98").Value = 2750
$ws.Range("M98").Value = 734.25
$ws.Range("N98").Value = -5746

$ws.Range("H107").Value = 7975.2144
$ws.Range("I107").Value = 7975.2144
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 7975.2144
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -6055.2144
$ws.Range("N107").ClearContents()

$ws.Range("H116").Value = 3440
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 4066.6667
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 4066.6667
$ws.Range("M116").Value = 942
$ws.Range("N116").Value = -10950.6667

$ws.Range("H122").Value = 1047.5
$ws.Range("I122").Value = 763.75
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 2291.25
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = 158.75
$ws.Range("N122").Value = -13150

$ws.Range("H132").Value = 11911947
$ws.Range("I132").Value = 6526
$ws.Range("J132").Value = 41675500
$ws.Range("K132").Value = 19578
$ws.Range("L132").Value = 125026500
$ws.Range("M132").Value = -17048
$ws.Range("N132").Value = -125031560

$ws.Range("H141").Value = 1464.375
$ws.Range("I141").Value = 1464.375
$ws.Range("K141").Value = 4393.125
$ws.Range("M141").Value = 786.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9692.610000000001
$ws.Range("I32").Value = 8779.337
$ws.Range("J32").Value = 17081.818
$ws.Range("K32").Value = 8779.337
$ws.Range("L32").Value = 17081.818
$ws.Range("M32").Value = -8492.337
$ws.Range("N32").Value = -17655.818

$ws.Range("H41").Value = 1524.5
$ws.Range("I41").Value = 1524.5
$ws.Range("K41").Value = 1524.5
$ws.Range("M41").Value = -1110.5

$ws.Range("H102").Value = 1832.6666
$ws.Range("I102").Value = 500
$ws.Range("J102").Value = 2499
$ws.Range("K102").Value = 500
$ws.Range("L102").Value = 2499
$ws.Range("M102").Value = 1122
$ws.Range("N102").Value = -5743

$ws.Range("H110").Value = 568.0476
$ws.Range("I110").Value = 568.0476
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 568.0476
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1476.9524
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 707.5517
$ws.Range("I94").Value = 504.7143
$ws.Range("J94").Value = 1240
$ws.Range("K94").Value = 504.7143
$ws.Range("L94").Value = 1240
$ws.Range("M94").Value = -53.71429999999998
$ws.Range("N94").Value = -2142

$ws.Range("H99").Value = 870.3333
$ws.Range("I99").Value = 800
$ws.Range("K99").Value = 800
$ws.Range("M99").Value = 698

$ws.Range("H134").Value = 5189.52
$ws.Range("I134").Value = 4555.1113
$ws.Range("J134").Value = 5546.375
$ws.Range("K134").Value = 13665.3339
$ws.Range("L134").Value = 16639.125
$ws.Range("M134").Value = -11130.3339
$ws.Range("N134").Value = -21709.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6551.4443
$ws.Range("I16").Value = 7829.5
$ws.Range("K16").Value = 7829.5
$ws.Range("M16").Value = -7542.5

$ws.Range("H99").Value = 2125
$ws.Range("I99").Value = 2333.3333
$ws.Range("K99").Value = 2333.3333
$ws.Range("M99").Value = -835.3332999999998

$ws.Range("H113").Value = 6551.4443
$ws.Range("I113").Value = 7829.5
$ws.Range("K113").Value = 7829.5
$ws.Range("M113").Value = -5659.5

$ws.Range("H123").Value = 29500
$ws.Range("J123").Value = 29500
$ws.Range("L123").Value = 29500
$ws.Range("N123").Value = -39300

$ws.Range("H126").Value = 2125
$ws.Range("I126").Value = 2333.3333
$ws.Range("K126").Value = 6999.999899999999
$ws.Range("M126").Value = -4529.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 490
$ws.Range("I20").Value = 490
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1470
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1243
$ws.Range("N20").ClearContents()

$ws.Range("H134").Value = 6258
$ws.Range("I134").Value = 3937
$ws.Range("J134").Value = 10900
$ws.Range("K134").Value = 11811
$ws.Range("L134").Value = 32700
$ws.Range("M134").Value = -6741
$ws.Range("N134").Value = -42840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 59326
$ws.Range("J125").Value = 59326
$ws.Range("L125").Value = 59326
$ws.Range("N125").Value = -64246

$ws.Range("H126").Value = 4544.9375
$ws.Range("I126").Value = 3353.6667
$ws.Range("J126").Value = 5259.7
$ws.Range("K126").Value = 10061.0001
$ws.Range("L126").Value = 15779.1
$ws.Range("M126").Value = -7591.000100000001
$ws.Range("N126").Value = -20719.1

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6899.048
$ws.Range("I7").Value = 8687.5
$ws.Range("J7").Value = 5798.4614
$ws.Range("K7").Value = 8687.5
$ws.Range("L7").Value = 5798.4614
$ws.Range("M7").Value = -8575.5
$ws.Range("N7").Value = -6022.4614

$ws.Range("H40").Value = 5456.696
$ws.Range("I40").Value = 7379.6
$ws.Range("J40").Value = 3977.5386
$ws.Range("K40").Value = 7379.6
$ws.Range("L40").Value = 3977.5386
$ws.Range("M40").Value = -7243.6
$ws.Range("N40").Value = -4249.5386

$ws.Range("H46").Value = 1137.3125
$ws.Range("I46").Value = 680
$ws.Range("J46").Value = 1289.75
$ws.Range("K46").Value = 680
$ws.Range("L46").Value = 1289.75
$ws.Range("M46").Value = -492
$ws.Range("N46").Value = -1665.75

$ws.Range("H126").Value = 6899.048
$ws.Range("I126").Value = 8687.5
$ws.Range("J126").Value = 5798.4614
$ws.Range("K126").Value = 26062.5
$ws.Range("L126").Value = 17395.3842
$ws.Range("M126").Value = -23592.5
$ws.Range("N126").Value = -22335.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2037.1333
$ws.Range("I126").Value = 2111.2144
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 6333.6432
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -3863.6432
$ws.Range("N126").Value = -7940

$ws.Range("H131").Value = 75999
$ws.Range("J131").Value = 75999
$ws.Range("L131").Value = 75999
$ws.Range("N131").Value = -86079
